# Daily attendance processing - normalize "Recorded By" (column G) ordering.
# For each entry in the comma-separated "Recorded By" list, any occurrence of
# the system account name ("System"/"system", case-insensitive) is moved to
# the end of the list (relative order of the system entries is reversed),
# while the other names keep their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") { continue }

    $parts = $val -split ',\s*'
    if ($parts.Count -le 1) { continue }

    $others = @()
    $systemParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $others += $p
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    # Manually reverse (the runtime's [array]::Reverse does not mutate in place here).
    $reversedSystemParts = @()
    for ($i = $systemParts.Count - 1; $i -ge 0; $i--) {
        $reversedSystemParts += $systemParts[$i]
    }

    $newParts = $others + $reversedSystemParts
    $newVal = [string]::Join(", ", $newParts)

    # NOTE: -eq/-ne (and even -ceq/-cne) on this runtime are case-insensitive,
    # so use the .NET string .Equals(...) method (case-sensitive, ordinal) to
    # decide whether a real change is needed.
    if (-not $val.Equals($newVal)) {
        $cell.Value = $newVal
    }
}
